# Desktop: added logs for 10 folder CV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("setseed")

# New "LSTM(256) dense(128)" row (14) values for the second (F:G) and third (J:K) result blocks
$ws.Range("F14").Value = "0.7260755181312561"
$ws.Range("G14").Value = "0.03157251168938171"
$ws.Range("J14").Value = "0.735030734539032"
$ws.Range("K14").Value = "0.03248890475003557"

# New "LSTM(256) dense(64)" row (15) values
$ws.Range("F15").Value = "0.7410886645317077"
$ws.Range("G15").Value = "0.035906858635031645"
$ws.Range("J15").Value = "0.726294994354248"
$ws.Range("K15").Value = "0.05806351534832032"

# New "LSTM(256) dense(32)" row (16) values
$ws.Range("F16").Value = "0.7290605902671814"
$ws.Range("G16").Value = "0.029760749967404104"
$ws.Range("J16").Value = "0.7203248381614685"
$ws.Range("K16").Value = "0.026293621346120167"

# New "LSTM(256,128)" row (17) values
$ws.Range("F17").Value = "0.7559262394905091"
$ws.Range("G17").Value = "0.012366525881863114"
$ws.Range("J17").Value = "0.7262071967124939"
$ws.Range("K17").Value = "0.051238683813622"

# Update the active selection on the sheet to match the saved cursor position
$ws.Activate()
$ws.Range("C17").Select()
